# Revert practice.xlsx stimuli paths back to their original (pre-"resources/"-prefix) form
# and restore the view to its original scroll position / selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count
$colCount = $usedRange.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $usedRange.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -is [string] -and $val.StartsWith("resources/stimuli/")) {
            $cell.Value2 = $val.Substring(10)
        }
    }
}

# Restore sheet view: scroll back to show column A (removes the saved topLeftCell="F1")
# and set the selection back to B2 (instead of K10).
[void]$ws.Activate()
[void]$ws.Range("A1").Select()
[void]$ws.Range("B2").Select()
